$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H107").Value = 225.7
$ws.Range("I107").Value = 224
$ws.Range("K107").Value = 224
$ws.Range("M107").Value = 1696

$ws.Range("H137").Value = 2238.0645
$ws.Range("I137").Value = 908.1667
$ws.Range("J137").Value = 4079.4614
$ws.Range("K137").Value = 2724.5001
$ws.Range("L137").Value = 12238.3842
$ws.Range("M137").Value = -174.5001000000002
$ws.Range("N137").Value = -17338.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8655.540999999999
$ws.Range("I32").Value = 6102.375
$ws.Range("K32").Value = 6102.375
$ws.Range("M32").Value = -5815.375

$ws.Range("H43").Value = 100000
$ws.Range("J43").Value = 100000
$ws.Range("L43").Value = 100000
$ws.Range("N43").Value = -100626

$ws.Range("H45").Value = 1483.1666
$ws.Range("I45").Value = 1116.3334
$ws.Range("J45").Value = 1850
$ws.Range("K45").Value = 1116.3334
$ws.Range("L45").Value = 1850
$ws.Range("M45").Value = -739.3334
$ws.Range("N45").Value = -2604

$ws.Range("H61").Value = 3094.1428
$ws.Range("I61").Value = 2942.5557
$ws.Range("J61").Value = 4003.6667
$ws.Range("K61").Value = 2942.5557
$ws.Range("L61").Value = 4003.6667
$ws.Range("M61").Value = -2730.5557
$ws.Range("N61").Value = -4427.6667

$ws.Range("H74").Value = 1475.4706
$ws.Range("I74").Value = 1135.7037
$ws.Range("K74").Value = 1135.7037
$ws.Range("M74").Value = -261.7037

$ws.Range("H77").Value = 1475.4706
$ws.Range("I77").Value = 1135.7037
$ws.Range("K77").Value = 5678.5185
$ws.Range("M77").Value = -1310.5185

$ws.Range("H88").Value = 1260
$ws.Range("I88").Value = 500
$ws.Range("J88").Value = 1450
$ws.Range("K88").Value = 500
$ws.Range("L88").Value = 1450
$ws.Range("M88").Value = -94
$ws.Range("N88").Value = -2262

$ws.Range("H91").Value = 1260
$ws.Range("I91").Value = 500
$ws.Range("J91").Value = 1450
$ws.Range("K91").Value = 500
$ws.Range("L91").Value = 1450
$ws.Range("M91").Value = 904
$ws.Range("N91").Value = -4258

$ws.Range("H122").Value = 716637.9399999999
$ws.Range("I122").Value = 1251678.9
$ws.Range("J122").Value = 3250
$ws.Range("K122").Value = 3755036.7
$ws.Range("L122").Value = 9750
$ws.Range("M122").Value = -3752586.7
$ws.Range("N122").Value = -14650

$ws.Range("H136").Value = 3094.1428
$ws.Range("I136").Value = 2942.5557
$ws.Range("J136").Value = 4003.6667
$ws.Range("K136").Value = 8827.667099999999
$ws.Range("L136").Value = 12011.0001
$ws.Range("M136").Value = -6277.667099999999
$ws.Range("N136").Value = -17111.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6069.5713
$ws.Range("I20").Value = 9124.25
$ws.Range("J20").Value = 1996.6666
$ws.Range("K20").Value = 9124.25
$ws.Range("L20").Value = 1996.6666
$ws.Range("M20").Value = -8877.25
$ws.Range("N20").Value = -2490.6666

$ws.Range("H86").Value = 3428.6365
$ws.Range("I86").Value = 3060.1428
$ws.Range("J86").Value = 4073.5
$ws.Range("K86").Value = 3060.1428
$ws.Range("L86").Value = 4073.5
$ws.Range("M86").Value = -1937.1428
$ws.Range("N86").Value = -6319.5

$ws.Range("H89").Value = 3428.6365
$ws.Range("I89").Value = 3060.1428
$ws.Range("J89").Value = 4073.5
$ws.Range("K89").Value = 15300.714
$ws.Range("L89").Value = 20367.5
$ws.Range("M89").Value = -9684.714
$ws.Range("N89").Value = -31599.5

$ws.Range("H134").Value = 2142.0527
$ws.Range("I134").Value = 1825.0625
$ws.Range("K134").Value = 5475.1875
$ws.Range("M134").Value = -2940.1875

$ws.Range("H141").Value = 49999
$ws.Range("I141").Value = 49999
$ws.Range("K141").Value = 49999
$ws.Range("M141").Value = -44819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 171.2
$ws.Range("J7").Value = 498.33334
$ws.Range("L7").Value = 498.33334
$ws.Range("N7").Value = -724.33334

$ws.Range("H31").Value = 3757.682
$ws.Range("I31").Value = 1372.5
$ws.Range("J31").Value = 5745.3335
$ws.Range("K31").Value = 1372.5
$ws.Range("L31").Value = 5745.3335
$ws.Range("M31").Value = -1077.5
$ws.Range("N31").Value = -6335.3335

$ws.Range("H34").Value = 3757.682
$ws.Range("I34").Value = 1372.5
$ws.Range("J34").Value = 5745.3335
$ws.Range("K34").Value = 1372.5
$ws.Range("L34").Value = 5745.3335
$ws.Range("M34").Value = -1170.5
$ws.Range("N34").Value = -6149.3335

$ws.Range("H86").Value = 7955.375
$ws.Range("I86").Value = 7173.5
$ws.Range("K86").Value = 7173.5
$ws.Range("M86").Value = -6050.5

$ws.Range("H89").Value = 7955.375
$ws.Range("I89").Value = 7173.5
$ws.Range("K89").Value = 35867.5
$ws.Range("M89").Value = -30251.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1018.7143
$ws.Range("I68").Value = 1274.5
$ws.Range("J68").Value = 916.4
$ws.Range("K68").Value = 3823.5
$ws.Range("L68").Value = 2749.2
$ws.Range("M68").Value = -3012.5
$ws.Range("N68").Value = -4371.2

$ws.Range("H70").Value = 3223.5
$ws.Range("I70").Value = 2964.6667
$ws.Range("K70").Value = 8894.000100000001
$ws.Range("M70").Value = -8579.000100000001

$ws.Range("H71").Value = 1018.7143
$ws.Range("I71").Value = 1274.5
$ws.Range("J71").Value = 916.4
$ws.Range("K71").Value = 11470.5
$ws.Range("L71").Value = 8247.6
$ws.Range("M71").Value = -7414.5
$ws.Range("N71").Value = -16359.6

$ws.Range("H73").Value = 3223.5
$ws.Range("I73").Value = 2964.6667
$ws.Range("K73").Value = 8894.000100000001
$ws.Range("M73").Value = -7802.000100000001

$ws.Range("H75").Value = 387.5
$ws.Range("J75").Value = 406.25
$ws.Range("L75").Value = 1218.75
$ws.Range("N75").Value = -3214.75

$ws.Range("H78").Value = 387.5
$ws.Range("J78").Value = 406.25
$ws.Range("L78").Value = 3656.25
$ws.Range("N78").Value = -13640.25

$ws.Range("H86").Value = 73.166664
$ws.Range("I86").Value = 67.25
$ws.Range("J86").Value = 85
$ws.Range("K86").Value = 201.75
$ws.Range("L86").Value = 255
$ws.Range("M86").Value = 984.25
$ws.Range("N86").Value = -2627

$ws.Range("H89").Value = 73.166664
$ws.Range("I89").Value = 67.25
$ws.Range("J89").Value = 85
$ws.Range("K89").Value = 605.25
$ws.Range("L89").Value = 765
$ws.Range("M89").Value = 5322.75
$ws.Range("N89").Value = -12621

$ws.Range("H107").Value = 341.26923
$ws.Range("J107").Value = 362.3913
$ws.Range("L107").Value = 1087.1739
$ws.Range("N107").Value = -4927.1739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6010.769
$ws.Range("I70").Value = 4890.8335
$ws.Range("J70").Value = 6970.7144
$ws.Range("K70").Value = 4890.8335
$ws.Range("L70").Value = 6970.7144
$ws.Range("M70").Value = -4620.8335
$ws.Range("N70").Value = -7510.7144

$ws.Range("H73").Value = 6010.769
$ws.Range("I73").Value = 4890.8335
$ws.Range("J73").Value = 6970.7144
$ws.Range("K73").Value = 4890.8335
$ws.Range("L73").Value = 6970.7144
$ws.Range("M73").Value = -3954.8335
$ws.Range("N73").Value = -8842.714400000001

$ws.Range("H97").Value = 681.8461
$ws.Range("J97").Value = 934.8333
$ws.Range("L97").Value = 934.8333
$ws.Range("N97").Value = -1926.8333

$ws.Range("H123").Value = 115714.43
$ws.Range("J123").Value = 115714.43
$ws.Range("L123").Value = 115714.43
$ws.Range("N123").Value = -120614.43

$ws.Range("H132").Value = 2454.7576
$ws.Range("I132").Value = 1661.8235
$ws.Range("J132").Value = 3297.25
$ws.Range("K132").Value = 4985.470499999999
$ws.Range("L132").Value = 9891.75
$ws.Range("M132").Value = -2455.470499999999
$ws.Range("N132").Value = -14951.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1086.5
$ws.Range("I22").Value = 743.625
$ws.Range("J22").Value = 1772.25
$ws.Range("K22").Value = 743.625
$ws.Range("L22").Value = 1772.25
$ws.Range("M22").Value = -448.625
$ws.Range("N22").Value = -2362.25

$ws.Range("H27").Value = 1086.5
$ws.Range("I27").Value = 743.625
$ws.Range("J27").Value = 1772.25
$ws.Range("K27").Value = 743.625
$ws.Range("L27").Value = 1772.25
$ws.Range("M27").Value = -636.625
$ws.Range("N27").Value = -1986.25

$ws.Range("H61").Value = 4389.5454
$ws.Range("I61").Value = 4389.5454
$ws.Range("K61").Value = 4389.5454
$ws.Range("M61").Value = -4187.5454

$ws.Range("H93").Value = 3074
$ws.Range("J93").Value = 3399
$ws.Range("L93").Value = 3399
$ws.Range("N93").Value = -5895

$ws.Range("H113").Value = 4389.5454
$ws.Range("I113").Value = 4389.5454
$ws.Range("K113").Value = 4389.5454
$ws.Range("M113").Value = -2219.5454

$ws.Range("H122").Value = 6280.5
$ws.Range("I122").Value = 6427.2856
$ws.Range("J122").Value = 6075
$ws.Range("K122").Value = 19281.8568
$ws.Range("L122").Value = 18225
$ws.Range("M122").Value = -16831.8568
$ws.Range("N122").Value = -23125

$ws.Range("H133").Value = 103749.75
$ws.Range("J133").Value = 103749.75
$ws.Range("L133").Value = 103749.75
$ws.Range("N133").Value = -108809.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 16232.6
$ws.Range("I45").Value = 15999
$ws.Range("J45").Value = 16291
$ws.Range("K45").Value = 15999
$ws.Range("L45").Value = 16291
$ws.Range("M45").Value = -15508
$ws.Range("N45").Value = -17273

$ws.Range("H122").Value = 3344.4285
$ws.Range("I122").Value = 3444
$ws.Range("K122").Value = 10332
$ws.Range("M122").Value = -7882

$ws.Range("H130").Value = 70993.8
$ws.Range("J130").Value = 70993.8
$ws.Range("L130").Value = 70993.8
$ws.Range("N130").Value = -81033.8
